$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 573.6667
$ws.Range("I33").Value = 684.7143
$ws.Range("J33").Value = 185
$ws.Range("K33").Value = 684.7143
$ws.Range("L33").Value = 185
$ws.Range("M33").Value = -455.7143
$ws.Range("N33").Value = -643

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 56032.4
$ws.Range("J40").Value = 65046.715
$ws.Range("L40").Value = 65046.715
$ws.Range("N40").Value = -65396.715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4233.1665
$ws.Range("I62").Value = 3800
$ws.Range("J62").Value = 4666.3335
$ws.Range("K62").Value = 3800
$ws.Range("L62").Value = 4666.3335
$ws.Range("M62").Value = -3176
$ws.Range("N62").Value = -5914.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4233.1665
$ws.Range("I65").Value = 3800
$ws.Range("J65").Value = 4666.3335
$ws.Range("K65").Value = 19000
$ws.Range("L65").Value = 23331.6675
$ws.Range("M65").Value = -15880
$ws.Range("N65").Value = -29571.6675

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 52631800
$ws.Range("I92").Value = 66666790
$ws.Range("K92").Value = 66666790
$ws.Range("M92").Value = -66665542

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 920.8125
$ws.Range("I107").Value = 916.4
$ws.Range("J107").Value = 987
$ws.Range("K107").Value = 916.4
$ws.Range("L107").Value = 987
$ws.Range("M107").Value = 1003.6
$ws.Range("N107").Value = -4827

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 44250.707
$ws.Range("I113").Value = 3413.5715
$ws.Range("J113").Value = 72836.7
$ws.Range("K113").Value = 3413.5715
$ws.Range("L113").Value = 72836.7
$ws.Range("M113").Value = -159.5715
$ws.Range("N113").Value = -79344.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14621.375
$ws.Range("I32").Value = 13829.233
$ws.Range("J32").Value = 26503.5
$ws.Range("K32").Value = 13829.233
$ws.Range("L32").Value = 26503.5
$ws.Range("M32").Value = -13542.233
$ws.Range("N32").Value = -27077.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3874.64
$ws.Range("I61").Value = 3052.4375
$ws.Range("K61").Value = 3052.4375
$ws.Range("M61").Value = -2840.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 107143660
$ws.Range("I74").Value = 150000660
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 150000660
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = -149999786
$ws.Range("N74").Value = -2948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 107143660
$ws.Range("I77").Value = 150000660
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 750003300
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = -749998932
$ws.Range("N77").Value = -14736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5957596.5
$ws.Range("I122").Value = 7250596
$ws.Range("J122").Value = 9799.6
$ws.Range("K122").Value = 21751788
$ws.Range("L122").Value = 29398.8
$ws.Range("M122").Value = -21749338
$ws.Range("N122").Value = -34298.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 71097
$ws.Range("J133").Value = 71097
$ws.Range("L133").Value = 71097
$ws.Range("N133").Value = -76157

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3874.64
$ws.Range("I136").Value = 3052.4375
$ws.Range("K136").Value = 9157.3125
$ws.Range("M136").Value = -6607.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 112871.5
$ws.Range("J140").Value = 112871.5
$ws.Range("L140").Value = 112871.5
$ws.Range("N140").Value = -123231.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 140000
$ws.Range("J141").Value = 140000
$ws.Range("L141").Value = 140000
$ws.Range("N141").Value = -150360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4037.45
$ws.Range("I105").Value = 4080.5
$ws.Range("K105").Value = 4080.5
$ws.Range("M105").Value = -2333.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 27779402
$ws.Range("I107").Value = 36112344
$ws.Range("J107").Value = 2926
$ws.Range("K107").Value = 36112344
$ws.Range("L107").Value = 2926
$ws.Range("M107").Value = -36110424
$ws.Range("N107").Value = -6766

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3117.1724
$ws.Range("I7").Value = 3089.3684
$ws.Range("J7").Value = 3170
$ws.Range("K7").Value = 3089.3684
$ws.Range("L7").Value = 3170
$ws.Range("M7").Value = -2976.3684
$ws.Range("N7").Value = -3396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14707840
$ws.Range("I31").Value = 17242726
$ws.Range("J31").Value = 5497.9
$ws.Range("K31").Value = 17242726
$ws.Range("L31").Value = 5497.9
$ws.Range("M31").Value = -17242431
$ws.Range("N31").Value = -6087.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14707840
$ws.Range("I34").Value = 17242726
$ws.Range("J34").Value = 5497.9
$ws.Range("K34").Value = 17242726
$ws.Range("L34").Value = 5497.9
$ws.Range("M34").Value = -17242524
$ws.Range("N34").Value = -5901.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 38975
$ws.Range("J97").Value = 38975
$ws.Range("L97").Value = 38975
$ws.Range("N97").Value = -40957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5950.5
$ws.Range("I99").Value = 2745.3333
$ws.Range("K99").Value = 2745.3333
$ws.Range("M99").Value = -1247.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2123.5715
$ws.Range("I122").Value = 1480.6333
$ws.Range("K122").Value = 4441.8999
$ws.Range("M122").Value = -1991.8999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5950.5
$ws.Range("I126").Value = 2745.3333
$ws.Range("K126").Value = 8235.999899999999
$ws.Range("M126").Value = -5765.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3020.6316
$ws.Range("I134").Value = 2500.1333
$ws.Range("K134").Value = 7500.3999
$ws.Range("M134").Value = -4965.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 102862.5
$ws.Range("J141").Value = 103603.125
$ws.Range("L141").Value = 103603.125
$ws.Range("N141").Value = -113963.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1713.2858
$ws.Range("J22").Value = 1482.1666
$ws.Range("L22").Value = 4446.4998
$ws.Range("N22").Value = -4784.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 1713.2858
$ws.Range("J27").Value = 1482.1666
$ws.Range("L27").Value = 4446.4998
$ws.Range("N27").Value = -4650.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 450
$ws.Range("I41").Value = 150
$ws.Range("J41").Value = 600
$ws.Range("K41").Value = 450
$ws.Range("L41").Value = 1800
$ws.Range("M41").Value = -112
$ws.Range("N41").Value = -2476

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 5706.4
$ws.Range("I134").Value = 1968.9231
$ws.Range("K134").Value = 5906.7693
$ws.Range("M134").Value = -836.7692999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3437.7
$ws.Range("I136").Value = 1841.7142
$ws.Range("J136").Value = 7161.6665
$ws.Range("K136").Value = 5525.142599999999
$ws.Range("L136").Value = 21484.9995
$ws.Range("M136").Value = -425.1425999999992
$ws.Range("N136").Value = -31684.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3360.5625
$ws.Range("I138").Value = 3360.5625
$ws.Range("K138").Value = 10081.6875
$ws.Range("M138").Value = -4941.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 337173.34
$ws.Range("I122").Value = 668429.6
$ws.Range("J122").Value = 5917.067
$ws.Range("K122").Value = 2005288.8
$ws.Range("L122").Value = 17751.201
$ws.Range("M122").Value = -2002838.8
$ws.Range("N122").Value = -22651.201

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3913.4814
$ws.Range("I132").Value = 3462.3333
$ws.Range("J132").Value = 5492.5
$ws.Range("K132").Value = 10386.9999
$ws.Range("L132").Value = 16477.5
$ws.Range("M132").Value = -7856.999899999999
$ws.Range("N132").Value = -21537.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1156.1428
$ws.Range("I93").Value = 1045.7368
$ws.Range("J93").Value = 1389.2222
$ws.Range("K93").Value = 1045.7368
$ws.Range("L93").Value = 1389.2222
$ws.Range("M93").Value = 202.2632000000001
$ws.Range("N93").Value = -3885.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7188.5713
$ws.Range("I122").Value = 3530.5454
$ws.Range("J122").Value = 9555.529
$ws.Range("K122").Value = 10591.6362
$ws.Range("L122").Value = 28666.587
$ws.Range("M122").Value = -8141.636200000001
$ws.Range("N122").Value = -33566.587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1024.6666
$ws.Range("J100").Value = 1858
$ws.Range("L100").Value = 3716
$ws.Range("N100").Value = -4798

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 701.5
$ws.Range("I113").Value = 658
$ws.Range("J113").Value = 797.2
$ws.Range("K113").Value = 1974
$ws.Range("L113").Value = 2391.6
$ws.Range("M113").Value = 196
$ws.Range("N113").Value = -6731.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3773.9429
$ws.Range("I122").Value = 1951.1305
$ws.Range("J122").Value = 7267.6665
$ws.Range("K122").Value = 5853.3915
$ws.Range("L122").Value = 21802.9995
$ws.Range("M122").Value = -3403.3915
$ws.Range("N122").Value = -26702.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1419.5333
$ws.Range("I132").Value = 1911.579
$ws.Range("J132").Value = 1191.5122
$ws.Range("K132").Value = 5734.737
$ws.Range("L132").Value = 3574.536599999999
$ws.Range("M132").Value = -3204.737
$ws.Range("N132").Value = -8634.536599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4162.8
$ws.Range("I136").Value = 4241.143
$ws.Range("J136").Value = 3980
$ws.Range("K136").Value = 12723.429
$ws.Range("L136").Value = 11940
$ws.Range("M136").Value = -10173.429
$ws.Range("N136").Value = -17040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 80020.8
$ws.Range("J140").Value = 80020.8
$ws.Range("L140").Value = 80020.8
$ws.Range("N140").Value = -90380.8
